$d = $word.ActiveDocument

# Locate the "Second option" paragraph by scanning paragraph text (robust
# against any paragraph re-numbering caused elsewhere in the document).
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Second option:*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Second option:' paragraph"
}

# --- Rewrite the "Second option" paragraph's text ---
$p = $paras.Item($targetIndex)
$start = $p.Range.Start
$end = $p.Range.End
$r = $d.Range($start, $end)
$r.Text = "Second option: Another menu inside Second option with choices -User should be able to 1.add, 2. delete files, 3. search file from the directory."

# --- Insert "4." at the very start of the following paragraph ---
# (re-fetch the Paragraphs collection since the text-length change above
# may have shifted character offsets)
$paras = $d.Paragraphs
$nextPara = $paras.Item($targetIndex + 1)
$insertionPoint = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)
$insertionPoint.InsertBefore("4.")
